$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.811.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.58%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.603.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.10%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'557.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.27%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'141.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.77%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.79%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.625.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.69%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.34%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.04%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +6.25%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.372"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +9.11%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'3.062.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +6.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'59.778.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.41%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.37%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.612.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.44%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'4.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.52%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'343.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.43%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +5.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +11.07%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.05%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.523"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +15.69%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'62.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.20%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.996"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.22%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.75%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0₃0785"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.12%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.08%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +1.94%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'6.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +2.46%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'158.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.30%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'19.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.16"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +3.95%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +5.54%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.920"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +5.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'37.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.66%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +1.98%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.848"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.67%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +2.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'294.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.41%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'140.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +12.85%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0981"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.602"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.42%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +4.23%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.65%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'10.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.21%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'4.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.61%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'19.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +3.00%  "
$ws.Range("E51").Style = "Normal"
